# srg-mapping-rhel9.xlsx: normalize the "IA Control" column (A) on the
# mapping sheet by de-duplicating the comma separated control IDs and
# joining them without the extra space after the comma. Also fills in the
# previously-empty Fix text for row 42 (M42) and rewords the Requirement
# text in F192 to reference "Red Hat Enterprise Linux 9" explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AU-4 (1),AU-4"
$ws.Range("A3").Value = "SC-5 (2),SC-5,CM-6 b"
$ws.Range("A4").Value = "AU-7 b,AU-12 (3),AU-8 b,CM-5 (1),AC-6 (9),AC-6 (8),AU-7 a"
$ws.Range("A5").Value = "AC-17 (9),CM-7 b,AC-17 (1),CM-6 b"
$ws.Range("A8").Value = "IA-2 (11),IA-2 (12)"
$ws.Range("A10").Value = "CM-7 (5) (b),CM-7 (2)"
$ws.Range("A12").Value = "AC-7 a,AC-7 b"
$ws.Range("A15").Value = "AU-3 (1),IA-8,IA-2"
$ws.Range("A16").Value = "AC-6 (10),CM-6 b"
$ws.Range("A17").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A19").Value = "IA-5 (1) (a),IA-5 (1) (b),CM-6 b"
$ws.Range("A21").Value = "MA-4 (7),SC-10,MA-4 e,AC-12"
$ws.Range("A22").Value = "AU-3,AU-14 (1),AU-12 a,AU-7 (1),CM-6 b,CM-5 (1),AU-6 (4),AU-3 (1),MA-4 (1) (a),AU-7 a"
$ws.Range("A25").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A29").Value = "SC-8 (2),SC-8 (1),SC-8"
$ws.Range("A31").Value = "AU-3,AU-12 a,AU-12 c,AC-2 (4),AU-3 (1),MA-4 (1) (a)"
$ws.Range("A34").Value = "AC-11 a,AC-11 b"
$ws.Range("A38").Value = "SI-11 b,AU-9"
$ws.Range("A39").Value = "AU-3,CM-6 b"
$ws.Range("A41").Value = "AU-4 (1)"
$ws.Range("A42").Value = "SC-28 (1),SC-28"
$ws.Range("A44").Value = "IA-11"
$ws.Range("A45").Value = "AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 b,AC-8 a"
$ws.Range("A48").Value = "CM-6 b"
$ws.Range("A50").Value = "IA-2 (5),CM-6 b"
$ws.Range("A53").Value = "SC-13,MA-4 (6)"
$ws.Range("A55").Value = "AC-17 (2),SC-8"
$ws.Range("A56").Value = "MA-4 (1) (a),AU-12 c"
$ws.Range("A63").Value = "AU-5 (1),AU-5 a"
$ws.Range("A65").Value = "IA-2 (2),CM-6 b"
$ws.Range("A67").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A69").Value = "AU-7 b,AU-12 (3),AU-12 a,AU-12 c,AU-8 b,CM-6 b,CM-5 (1),AU-7 a"
$ws.Range("A71").Value = "AU-3,AU-4 (1)"
$ws.Range("A77").Value = "AU-3,AU-12 a,AU-12 c,AC-2 (4),AU-3 (1),MA-4 (1) (a)"
$ws.Range("A79").Value = "AU-9 (3),AU-9"
$ws.Range("A80").Value = "IA-2 (1),IA-2 (4),IA-2 (2),IA-2 (3)"
$ws.Range("A81").Value = "CM-5 (3),CM-6 b"
$ws.Range("A86").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A88").Value = "CM-5 (1),AC-2 (4),AC-6 (9),AU-12 c"
$ws.Range("A89").Value = "IA-2 (4),IA-2 (5),IA-2 (2),IA-2 (3),IA-2"
$ws.Range("A90").Value = "IA-2 (11),IA-2 (12)"
$ws.Range("A91").Value = "AU-9 (3),AU-9"
$ws.Range("A96").Value = "AC-18 (1),SC-8 (1),SC-8"
$ws.Range("A97").Value = "AU-8 b,AU-8 (1) (a),AU-8 (1) (b)"
$ws.Range("A99").Value = "AU-9"
$ws.Range("A101").Value = "IA-11,AC-3 (4)"
$ws.Range("A102").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A111").Value = "AU-5 b,AU-5 a"
$ws.Range("A119").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A123").Value = "CM-7 b,CM-7 a"
$ws.Range("A124").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A125").Value = "AC-18 (1),CM-7 a"
$ws.Range("A128").Value = "CM-6 b,IA-5 (1) (c),CM-7 a"
$ws.Range("A136").Value = "AC-11 (1),AC-11 b"
$ws.Range("A139").Value = "SI-6 d,SI-6 b,CM-3 (5)"
$ws.Range("A142").Value = "AC-2 (2)"
$ws.Range("A148").Value = "AU-3,AU-14 (1),AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A157").Value = "AU-3,AU-12 a,AU-12 c,AU-3 (1),MA-4 (1) (a)"
$ws.Range("A159").Value = "AC-17 (2),SC-8"
$ws.Range("A175").Value = "SI-16,CM-7 a"
$ws.Range("A181").Value = "SC-3,CM-6 b"
$ws.Range("F192").Value = "Red Hat Enterprise Linux 9 must protect the confidentiality and integrity of all information at rest."

$ws.Range("M42").Value = "Configure Red Hat Enterprise Linux 9 to prevent unauthorized modification of all information at rest by using disk encryption.`nEncrypting a partition in an already installed system is more difficult, because existing partitions will need to be resized and changed.`nTo encrypt an entire partition, dedicate a partition for encryption in the partition layout."
